$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 40
$ws.Range("B1").Value = 40
$ws.Range("B3").Value = 185
$ws.Range("E4").Value = 3

$ws.Activate()
$ws.Range("E3").Select()
